$wb = $excel.ActiveWorkbook

# --- "Logs" sheet: append the new mail-log entry as row 10 ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A10").Value = "Terugbetaling"
$logs.Range("B10").Value = "support@testbedrijf123.nl"
$logs.Range("D10").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("F10").Value = "2025-08-19 19:47:54"
$logs.Range("G10").Value = "Nee"
$logs.Range("H10").Value = "Ja"
$logs.Range("I10").Value = "Nee"
$logs.Range("J10").Value = "Nee"

# --- Extend the conditional-formatting ranges (D/G/H/I/J) down to row 10 ---
$columns = @("D", "G", "H", "I", "J")
foreach ($col in $columns) {
    $oldRange = $logs.Range($col + "2:" + $col + "9")
    $newRange = $logs.Range($col + "2:" + $col + "10")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- "Dashboard" sheet: bump the summary count for this category ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 9
